# Agenda: generacion de puntos.
#
# The SharePoint "content type" custom XML parts that ship inside this
# template got re-synced, which swaps which physical item slot
# (customXml/item1.xml <-> customXml/item2.xml, and their matching
# itemProps1.xml <-> itemProps2.xml companions) holds which payload:
#   - the "FormTemplates" content-type forms stub
#     (http://schemas.microsoft.com/sharepoint/v3/contenttype/forms)
#   - the document-management "properties" payload (TaxCatchAll /
#     lcf76f155ced4ddcb4097134ff3c332f taxonomy field)
#
# Use the CustomXMLParts object model (the supported way to touch
# customXml/*) to relocate each payload into the other's slot instead of
# touching the package XML directly.

$p = $ppt.ActivePresentation
$parts = $p.CustomXMLParts

$formsNs = "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
$propsNs = "http://schemas.microsoft.com/office/2006/metadata/properties"

$formsPart = $parts.SelectByNamespace($formsNs).Item(1)
$propsPart = $parts.SelectByNamespace($propsNs).Item(1)

$formsXml = $formsPart.XML
$propsXml = $propsPart.XML

# Swap the payloads between the two parts/slots.
$formsPart.XML = $propsXml
$propsPart.XML = $formsXml
